$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.761.47"
$ws.Range("E2").Value = "  -1.81%  "
$ws.Range("D3").Value = "3.557.56"
$ws.Range("E3").Value = "  -2.80%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.94%  "
$ws.Range("D7").Value = "3.554.11"
$ws.Range("E7").Value = "  -2.72%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.617"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.68%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("E10").Value = "  -1.04%  "
$ws.Range("E11").Value = "  -4.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.55"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.80%  "
$ws.Range("E13").Value = "  +1.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.23%  "
$ws.Range("D15").Value = "4.131.00"
$ws.Range("E15").Value = "  -2.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.03%  "
$ws.Range("D17").Value = "3.555.71"
$ws.Range("E17").Value = "  -2.79%  "
$ws.Range("D18").Value = "69.757.57"
$ws.Range("E18").Value = "  -1.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.75%  "
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("E21").Value = "  -4.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "493.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "95.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.07%  "
$ws.Range("E28").Value = "  -7.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "66.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.114"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "564.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -10.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "38.45"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.69%  "
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0783"
$ws.Range("E39").Value = "  -6.04%  "
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.391"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.80%  "
$ws.Range("E41").Value = "  -2.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.13"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.09%  "
$ws.Range("E43").Value = "  -10.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.13%  "
$ws.Range("D45").Value = "3.218.66"
$ws.Range("E45").Value = "  -3.12%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.21%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0438"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("E49").Value = "  -3.23%  "
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("E51").Value = "  -5.04%  "
